# Daily attendance processing - reorder the "Recorded By" (column G) names.
#
# For every row in the "Session Analysis Results" sheet whose column G
# contains a comma-separated list of recorders (e.g. "user@mail.com, System"),
# flip the order of the names (e.g. "System, user@mail.com"). Rows with a
# single recorder are left untouched (reversing one item is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$lastRow   = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val  = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Count -gt 1) {
            $reversedParts = $parts[($parts.Count - 1)..0]
            $newVal = [string]::Join(", ", $reversedParts)
            $cell.Value2 = $newVal
        }
    }
}
